$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data cells with the refreshed OR import dataset
$ws.Cells.Item(1,1).Value = "group"
$ws.Cells.Item(1,2).Value = "OR"
$ws.Cells.Item(1,3).Value = "lower"
$ws.Cells.Item(1,4).Value = "higher"
$ws.Cells.Item(1,5).Value = "Adj"
$ws.Cells.Item(2,1).Value = "Age in 10-year scale"
$ws.Cells.Item(2,2).Value = 1.04
$ws.Cells.Item(2,3).Value = 0.99
$ws.Cells.Item(2,4).Value = 1.09
$ws.Cells.Item(2,5).Value = "Unadjusted"
$ws.Cells.Item(3,1).Value = "Sex: Male (ref: Female)"
$ws.Cells.Item(3,2).Value = 1.12
$ws.Cells.Item(3,3).Value = 1.02
$ws.Cells.Item(3,4).Value = 1.23
$ws.Cells.Item(3,5).Value = "Unadjusted"
$ws.Cells.Item(4,1).Value = "U/R: Urban Core (ref: Rural)"
$ws.Cells.Item(4,2).Value = 0.87
$ws.Cells.Item(4,3).Value = 0.72
$ws.Cells.Item(4,4).Value = 1.05
$ws.Cells.Item(4,5).Value = "Unadjusted"
$ws.Cells.Item(5,1).Value = "U/R: Urban Fringe (ref: Rural)"
$ws.Cells.Item(5,2).Value = 0.88
$ws.Cells.Item(5,3).Value = 0.61
$ws.Cells.Item(5,4).Value = 1.26
$ws.Cells.Item(5,5).Value = "Unadjusted"
$ws.Cells.Item(6,1).Value = "U/R: Urban out CMA/CAs (ref: Rural)"
$ws.Cells.Item(6,2).Value = 1.49
$ws.Cells.Item(6,3).Value = 0.8
$ws.Cells.Item(6,4).Value = 2.86
$ws.Cells.Item(6,5).Value = "Unadjusted"
$ws.Cells.Item(7,1).Value = "U/R: 2nd Urban Core (ref: Rural)"
$ws.Cells.Item(7,2).Value = 0.38
$ws.Cells.Item(7,3).Value = 0.25
$ws.Cells.Item(7,4).Value = 0.58
$ws.Cells.Item(7,5).Value = "Unadjusted"
$ws.Cells.Item(8,1).Value = "U/R: DA Only (ref: Rural)"
$ws.Cells.Item(8,2).Value = 0.8
$ws.Cells.Item(8,3).Value = 0.53
$ws.Cells.Item(8,4).Value = 1.21
$ws.Cells.Item(8,5).Value = "Unadjusted"
$ws.Cells.Item(9,1).Value = "Non-White (ref: White)"
$ws.Cells.Item(9,2).Value = 0.63
$ws.Cells.Item(9,3).Value = 0.49
$ws.Cells.Item(9,4).Value = 0.8
$ws.Cells.Item(9,5).Value = "Unadjusted"
$ws.Cells.Item(10,1).Value = "Edu: High School (ref: <High School)"
$ws.Cells.Item(10,2).Value = 1.38
$ws.Cells.Item(10,3).Value = 1.06
$ws.Cells.Item(10,4).Value = 1.81
$ws.Cells.Item(10,5).Value = "Unadjusted"
$ws.Cells.Item(11,1).Value = "Edu: Vocational Tr (ref: <High School)"
$ws.Cells.Item(11,2).Value = 1.1
$ws.Cells.Item(11,3).Value = 0.82
$ws.Cells.Item(11,4).Value = 1.47
$ws.Cells.Item(11,5).Value = "Unadjusted"
$ws.Cells.Item(12,1).Value = "Edu: Non-Uni Cert (ref: <High School)"
$ws.Cells.Item(12,2).Value = 0.91
$ws.Cells.Item(12,3).Value = 0.7
$ws.Cells.Item(12,4).Value = 1.19
$ws.Cells.Item(12,5).Value = "Unadjusted"
$ws.Cells.Item(13,1).Value = "Edu: Bachelor (ref: <High School)"
$ws.Cells.Item(13,2).Value = 1.35
$ws.Cells.Item(13,3).Value = 1.05
$ws.Cells.Item(13,4).Value = 1.74
$ws.Cells.Item(13,5).Value = "Unadjusted"
$ws.Cells.Item(14,1).Value = "Edu: Graduate (ref: <High School)"
$ws.Cells.Item(14,2).Value = 1.27
$ws.Cells.Item(14,3).Value = 0.98
$ws.Cells.Item(14,4).Value = 1.64
$ws.Cells.Item(14,5).Value = "Unadjusted"
$ws.Cells.Item(15,1).Value = "Depression Scale"
$ws.Cells.Item(15,2).Value = 1
$ws.Cells.Item(15,3).Value = 1
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(15,5).Value = "Unadjusted"
$ws.Cells.Item(16,1).Value = "Traveling Distance (km)"
$ws.Cells.Item(16,2).Value = 0.9734
$ws.Cells.Item(16,3).Value = 0.97
$ws.Cells.Item(16,4).Value = 0.98
$ws.Cells.Item(16,5).Value = "Unadjusted"
$ws.Cells.Item(17,1).Value = "Vaccination >15 Days"
$ws.Cells.Item(17,2).Value = 0.62
$ws.Cells.Item(17,3).Value = 0.57
$ws.Cells.Item(17,4).Value = 0.68
$ws.Cells.Item(17,5).Value = "Unadjusted"
$ws.Cells.Item(18,1).Value = "Prov Group: NF & NS (ref: BC)"
$ws.Cells.Item(18,2).Value = 1.63
$ws.Cells.Item(18,3).Value = 1.4
$ws.Cells.Item(18,4).Value = 1.9
$ws.Cells.Item(18,5).Value = "Unadjusted"
$ws.Cells.Item(19,1).Value = "Prov Group: QC (ref: BC)"
$ws.Cells.Item(19,2).Value = 0.72
$ws.Cells.Item(19,3).Value = 0.62
$ws.Cells.Item(19,4).Value = 0.84
$ws.Cells.Item(19,5).Value = "Unadjusted"
$ws.Cells.Item(20,1).Value = "Prov Group: ON (ref: BC)"
$ws.Cells.Item(20,2).Value = 2.14
$ws.Cells.Item(20,3).Value = 1.85
$ws.Cells.Item(20,4).Value = 2.47
$ws.Cells.Item(20,5).Value = "Unadjusted"
$ws.Cells.Item(21,1).Value = "Prov Group: MT & AB (ref: BC)"
$ws.Cells.Item(21,2).Value = 1.94
$ws.Cells.Item(21,3).Value = 1.67
$ws.Cells.Item(21,4).Value = 2.25
$ws.Cells.Item(21,5).Value = "Unadjusted"
$ws.Cells.Item(22,1).Value = "Outbreak Peak"
$ws.Cells.Item(22,2).Value = 2.18
$ws.Cells.Item(22,3).Value = 1.75
$ws.Cells.Item(22,4).Value = 2.74
$ws.Cells.Item(22,5).Value = "Unadjusted"
$ws.Cells.Item(23,1).Value = "Age in 10-year scale"
$ws.Cells.Item(23,2).Value = 1.04
$ws.Cells.Item(23,3).Value = 0.98
$ws.Cells.Item(23,4).Value = 1.1
$ws.Cells.Item(23,5).Value = "Adjusted"
$ws.Cells.Item(24,1).Value = "Sex: Male (ref: Female)"
$ws.Cells.Item(24,2).Value = 1.1
$ws.Cells.Item(24,3).Value = 1
$ws.Cells.Item(24,4).Value = 1.21
$ws.Cells.Item(24,5).Value = "Adjusted"
$ws.Cells.Item(25,1).Value = "U/R: Urban Core (ref: Rural)"
$ws.Cells.Item(25,2).Value = 0.49
$ws.Cells.Item(25,3).Value = 0.4
$ws.Cells.Item(25,4).Value = 0.61
$ws.Cells.Item(25,5).Value = "Adjusted"
$ws.Cells.Item(26,1).Value = "U/R: Urban Fringe (ref: Rural)"
$ws.Cells.Item(26,2).Value = 0.66
$ws.Cells.Item(26,3).Value = 0.45
$ws.Cells.Item(26,4).Value = 0.97
$ws.Cells.Item(26,5).Value = "Adjusted"
$ws.Cells.Item(27,1).Value = "U/R: Urban out CMA/CAs (ref: Rural)"
$ws.Cells.Item(27,2).Value = 1.59
$ws.Cells.Item(27,3).Value = 0.82
$ws.Cells.Item(27,4).Value = 3.17
$ws.Cells.Item(27,5).Value = "Adjusted"
$ws.Cells.Item(28,1).Value = "U/R: 2nd Urban Core (ref: Rural)"
$ws.Cells.Item(28,2).Value = 0.51
$ws.Cells.Item(28,3).Value = 0.32
$ws.Cells.Item(28,4).Value = 0.79
$ws.Cells.Item(28,5).Value = "Adjusted"
$ws.Cells.Item(29,1).Value = "U/R: DA Only (ref: Rural)"
$ws.Cells.Item(29,2).Value = 0.47
$ws.Cells.Item(29,3).Value = 0.3
$ws.Cells.Item(29,4).Value = 0.73
$ws.Cells.Item(29,5).Value = "Adjusted"
$ws.Cells.Item(30,1).Value = "Non-White (ref: White)"
$ws.Cells.Item(30,2).Value = 0.76
$ws.Cells.Item(30,3).Value = 0.57
$ws.Cells.Item(30,4).Value = 1
$ws.Cells.Item(30,5).Value = "Adjusted"
$ws.Cells.Item(31,1).Value = "Edu: High School (ref: <High School)"
$ws.Cells.Item(31,2).Value = 1.37
$ws.Cells.Item(31,3).Value = 1.03
$ws.Cells.Item(31,4).Value = 1.83
$ws.Cells.Item(31,5).Value = "Adjusted"
$ws.Cells.Item(32,1).Value = "Edu: Vocational Tr (ref: <High School)"
$ws.Cells.Item(32,2).Value = 1.11
$ws.Cells.Item(32,3).Value = 0.81
$ws.Cells.Item(32,4).Value = 1.52
$ws.Cells.Item(32,5).Value = "Adjusted"
$ws.Cells.Item(33,1).Value = "Edu: Non-Uni Cert (ref: <High School)"
$ws.Cells.Item(33,2).Value = 0.91
$ws.Cells.Item(33,3).Value = 0.69
$ws.Cells.Item(33,4).Value = 1.21
$ws.Cells.Item(33,5).Value = "Adjusted"
$ws.Cells.Item(34,1).Value = "Edu: Bachelor (ref: <High School)"
$ws.Cells.Item(34,2).Value = 1.26
$ws.Cells.Item(34,3).Value = 0.96
$ws.Cells.Item(34,4).Value = 1.66
$ws.Cells.Item(34,5).Value = "Adjusted"
$ws.Cells.Item(35,1).Value = "Edu: Graduate (ref: <High School)"
$ws.Cells.Item(35,2).Value = 1.16
$ws.Cells.Item(35,3).Value = 0.88
$ws.Cells.Item(35,4).Value = 1.52
$ws.Cells.Item(35,5).Value = "Adjusted"
$ws.Cells.Item(36,1).Value = "Depression Scale"
$ws.Cells.Item(36,2).Value = 1
$ws.Cells.Item(36,3).Value = 1
$ws.Cells.Item(36,4).Value = 1
$ws.Cells.Item(36,5).Value = "Adjusted"
$ws.Cells.Item(37,1).Value = "Traveling Distance (km)"
$ws.Cells.Item(37,2).Value = 0.96
$ws.Cells.Item(37,3).Value = 0.96
$ws.Cells.Item(37,4).Value = 0.97
$ws.Cells.Item(37,5).Value = "Adjusted"
$ws.Cells.Item(38,1).Value = "Vaccination >15 Days"
$ws.Cells.Item(38,2).Value = 0.64
$ws.Cells.Item(38,3).Value = 0.58
$ws.Cells.Item(38,4).Value = 0.71
$ws.Cells.Item(38,5).Value = "Adjusted"
$ws.Cells.Item(39,1).Value = "Prov Group: NF & NS (ref: BC)"
$ws.Cells.Item(39,2).Value = 1.48
$ws.Cells.Item(39,3).Value = 1.26
$ws.Cells.Item(39,4).Value = 1.73
$ws.Cells.Item(39,5).Value = "Adjusted"
$ws.Cells.Item(40,1).Value = "Prov Group: QC (ref: BC)"
$ws.Cells.Item(40,2).Value = 0.77
$ws.Cells.Item(40,3).Value = 0.65
$ws.Cells.Item(40,4).Value = 0.9
$ws.Cells.Item(40,5).Value = "Adjusted"
$ws.Cells.Item(41,1).Value = "Prov Group: ON (ref: BC)"
$ws.Cells.Item(41,2).Value = 2.52
$ws.Cells.Item(41,3).Value = 2.16
$ws.Cells.Item(41,4).Value = 2.94
$ws.Cells.Item(41,5).Value = "Adjusted"
$ws.Cells.Item(42,1).Value = "Prov Group: MT & AB (ref: BC)"
$ws.Cells.Item(42,2).Value = 2.04
$ws.Cells.Item(42,3).Value = 1.74
$ws.Cells.Item(42,4).Value = 2.38
$ws.Cells.Item(42,5).Value = "Adjusted"
$ws.Cells.Item(43,1).Value = "Outbreak Peak"
$ws.Cells.Item(43,2).Value = 1.73
$ws.Cells.Item(43,3).Value = 1.36
$ws.Cells.Item(43,4).Value = 2.2
$ws.Cells.Item(43,5).Value = "Adjusted"

# Widen column A to fit the longer relabeled group names
$ws.Range("A:A").ColumnWidth = 39.166666666666664

# Update the active selection (also resets scroll/top-left cell to default)
$ws.Range("A29").Select()
